$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Row 2 (Beta) Hp(0.07)/D_eye/Hp(10) Mean/HDI columns ---
$ws.Range("F2").Value = 11.83282077082337
$ws.Range("G2").Value = 11.50355499617784
$ws.Range("H2").Value = 12.15969956504401
$ws.Range("I2").Value = 1.958846661237457
$ws.Range("J2").Value = 1.9376149628667
$ws.Range("K2").Value = 1.979667720669618
$ws.Range("L2").Value = 0.152379779718232
$ws.Range("M2").Value = 0.1507092671500617
$ws.Range("N2").Value = 0.1540427797154021

# --- Update existing Row 3 (Gamma) Hp(0.07)/D_eye/Hp(10) Mean/HDI columns ---
$ws.Range("F3").Value = 0.001988920169596953
$ws.Range("G3").Value = 0.001225239106697428
$ws.Range("H3").Value = 0.002928094946145958
$ws.Range("I3").Value = 0.001842583827726434
$ws.Range("J3").Value = 0.001127358061955948
$ws.Range("K3").Value = 0.002719638412528947
$ws.Range("L3").Value = 0.002073003149112929
$ws.Range("M3").Value = 0.001295068847266368
$ws.Range("N3").Value = 0.003025949016168607

# --- Add new Row 4 (Beta + Gamma) ---
$ws.Range("A4").Value = 2
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 12.00687180793019
$ws.Range("D4").Value = 1.974477778970852
$ws.Range("E4").Value = 0.1537386519519979
$ws.Range("F4").Value = 11.83480969099296
$ws.Range("G4").Value = 11.50478023528454
$ws.Range("H4").Value = 12.16262765999016
$ws.Range("I4").Value = 1.960689245065184
$ws.Range("J4").Value = 1.938742320928656
$ws.Range("K4").Value = 1.982387359082146
$ws.Range("L4").Value = 0.1544527828673449
$ws.Range("M4").Value = 0.1520043359973281
$ws.Range("N4").Value = 0.1570687287315708
